$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experiment 1")

# xlPasteFormats = -4122
$xlPasteFormats = -4122

# --- Build the new "Experiment # X" / IMCS Compression matrix under the
# --- existing tables, in columns E:G, starting at row 13 (mirrors the
# --- layout already used for "Experiment # 1" / "Experiment # 2" above it).
# Formats are copied one cell at a time (rather than as a 3-wide block) so
# the destination reuses the exact same style records as the source band
# instead of the engine minting new (merged-border) style variants.

function Copy-CellFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
}

# Merge the title band first (like E1:G1 above it) so the destination
# already has the same merged shape as the source when formats are pasted.
$ws.Range("E13:G13").Merge()

# Row 13 -> title band (copy format from the row-1 title band)
Copy-CellFormat "E1" "E13"
Copy-CellFormat "F1" "F13"
Copy-CellFormat "G1" "G13"

# Row 14 -> header band (copy format from the row-2 header band)
Copy-CellFormat "E2" "E14"
Copy-CellFormat "F2" "F14"
Copy-CellFormat "G2" "G14"

# Rows 15-17 -> first 3-row striped group (copy format from rows 3-5)
Copy-CellFormat "E3" "E15"
Copy-CellFormat "F3" "F15"
Copy-CellFormat "G3" "G15"
Copy-CellFormat "E4" "E16"
Copy-CellFormat "F4" "F16"
Copy-CellFormat "G4" "G16"
Copy-CellFormat "E5" "E17"
Copy-CellFormat "F5" "F17"
Copy-CellFormat "G5" "G17"

# Rows 18-20 -> second 3-row striped group (copy format from rows 3-5 again)
Copy-CellFormat "E3" "E18"
Copy-CellFormat "F3" "F18"
Copy-CellFormat "G3" "G18"
Copy-CellFormat "E4" "E19"
Copy-CellFormat "F4" "F19"
Copy-CellFormat "G4" "G19"
Copy-CellFormat "E5" "E20"
Copy-CellFormat "F5" "F20"
Copy-CellFormat "G5" "G20"

$excel.CutCopyMode = $false

# --- Title row ---
$ws.Range("E13").Value = "Experiment # X"

# --- Row numbers (E column) and Split header ---
$ws.Range("E14").Value = "Split"
$ws.Range("E15").Value = 1
$ws.Range("E16").Value = 2
$ws.Range("E17").Value = 3
$ws.Range("E18").Value = 4
$ws.Range("E19").Value = 5
$ws.Range("E20").Value = 6

# --- IMCS Compression column (values first, header after - matches
# --- original authoring order captured by the shared-string table) ---
$ws.Range("F15").Value = "NO MEMCOMPRESS"
$ws.Range("F16").Value = "MEMCOMPRESS FOR DML"
$ws.Range("F17").Value = "MEMCOMPRESS FOR QUERY LOW"
$ws.Range("F14").Value = "IMCS Compression"

# --- INMEMORY_SIZE column ---
$ws.Range("G14").Value = "INMEMORY_SIZE"
$ws.Range("G15").Value = "0G"

$ws.Range("F18").Value = "MEMCOMPRESS FOR QUERY HIGH"
$ws.Range("F19").Value = "MEMCOMPRESS FOR CAPACITY LOW"
$ws.Range("F20").Value = "MEMCOMPRESS FOR CAPACITY HIGH"

$ws.Range("G16").Value = "64G"
$ws.Range("G17").Value = "64G"
$ws.Range("G18").Value = "64G"
$ws.Range("G19").Value = "64G"
$ws.Range("G20").Value = "64G"

# --- Selection / scroll state to mirror the authored workbook ---
$ws.Activate()
$ws.Range("G23").Select()
